$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.040.95'
$ws.Range("E2").Value = '  -1.69%  '
$ws.Range("D3").Value = '1.791.17'
$ws.Range("E3").Value = '  -0.72%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''223.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").Value = '''0.549'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.74%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").Value = '''32.20'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.09%  '
$ws.Range("D9").Value = '''0.283'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.59%  '
$ws.Range("D10").Value = '''0.0706'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.53%  '
$ws.Range("E11").Value = '  -0.01%  '
$ws.Range("D12").Value = '2.049.52'
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("D13").Value = '1.790.85'
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("D14").Value = '''10.79'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.93%  '
$ws.Range("D15").Value = '''0.623'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.18%  '
$ws.Range("D16").Value = '34.049.96'
$ws.Range("E16").Value = '  -1.69%  '
$ws.Range("D17").Value = '''4.14'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.32%  '
$ws.Range("D18").Value = '''67.79'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.13%  '
$ws.Range("D19").Value = '''242.70'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = '0.0₃0780'
$ws.Range("E20").Value = '  -2.60%  '
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("D22").Value = '''10.64'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.04%  '
$ws.Range("D23").Value = '''4.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.58%  '
$ws.Range("E24").Value = '  -2.85%  '
$ws.Range("D25").Value = '''158.29'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.02%  '
$ws.Range("D26").Value = '''16.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.23%  '
$ws.Range("D27").Value = '''6.98'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.45%  '
$ws.Range("E28").Value = '  -2.28%  '
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("D30").Value = '''0.0517'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.27%  '
$ws.Range("E31").Value = '  +0.16%  '
$ws.Range("D32").Value = '''3.65'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.09%  '
$ws.Range("D33").Value = '''3.48'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.28%  '
$ws.Range("D34").Value = '''1.80'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.87%  '
$ws.Range("D35").Value = '1.382.78'
$ws.Range("E35").Value = '  -3.66%  '
$ws.Range("D36").Value = '''0.643'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("E37").Value = '  -2.08%  '
$ws.Range("E38").Value = '  -4.30%  '
$ws.Range("B39").Value = 'HuobiToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D39").Value = '''2.35'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.20%  '
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").Value = '''79.21'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.70%  '
$ws.Range("E41").Value = '  -3.68%  '
$ws.Range("D42").Value = '''0.910'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.88%  '
$ws.Range("D43").Value = '''2.15'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.33%  '
$ws.Range("D44").Value = '0.0₆0136'
$ws.Range("E44").Value = '  +7.39%  '
$ws.Range("D45").Value = '''0.0497'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.26%  '
$ws.Range("E46").Value = '  -0.72%  '
$ws.Range("D47").Value = '''106.99'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '1.948.18'
$ws.Range("E48").Value = '  -0.50%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '''5.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.57%  '
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("D51").Value = '''11.91'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.97%  '
